$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.01139766666666667
$ws.Range("H2").Value = 0.034193
$ws.Range("I2").Value = 0.1481191086775714
$ws.Range("J2").Value = 0.1481191086775714
$ws.Range("M2").Value = 9.162874
$ws.Range("N2").Value = 27.488622
$ws.Range("O2").Value = 0.3895809024642617
$ws.Range("P2").Value = 0.3895809024642617
$ws.Range("Q2").Value = 0.1044353835606667
$ws.Range("R2").Value = 0.939918452046
$ws.Range("S2").Value = 0.05770437603081031
$ws.Range("T2").Value = 0.05770437603081031
$ws.Range("G3").Value = 0.01139766666666667
$ws.Range("H3").Value = 0.034193
$ws.Range("I3").Value = 0.1481191086775714
$ws.Range("J3").Value = 0.1481191086775714
$ws.Range("M3").Value = 8.050406666666667
$ws.Range("O3").Value = 0.3422817660053286
$ws.Range("P3").Value = 0.3422817660053286
$ws.Range("Q3").Value = 0.0917558517177778
$ws.Range("R3").Value = 0.8258026654600001
$ws.Range("S3").Value = 0.05069847009729432
$ws.Range("T3").Value = 0.05069847009729432
$ws.Range("G4").Value = 0.01139766666666667
$ws.Range("H4").Value = 0.034193
$ws.Range("I4").Value = 0.1481191086775714
$ws.Range("J4").Value = 0.1481191086775714
$ws.Range("M4").Value = 6.306542666666666
$ws.Range("N4").Value = 18.919628
$ws.Range("O4").Value = 0.2681373315304097
$ws.Range("P4").Value = 0.2681373315304097
$ws.Range("Q4").Value = 0.07187987113377778
$ws.Range("R4").Value = 0.646918840204
$ws.Range("S4").Value = 0.03971626254946675
$ws.Range("T4").Value = 0.03971626254946675
$ws.Range("G5").Value = 0.015206
$ws.Range("H5").Value = 0.045618
$ws.Range("I5").Value = 0.197610548932631
$ws.Range("J5").Value = 0.197610548932631
$ws.Range("M5").Value = 9.162874
$ws.Range("N5").Value = 27.488622
$ws.Range("O5").Value = 0.3895809024642617
$ws.Range("P5").Value = 0.3895809024642617
$ws.Range("Q5").Value = 0.139330662044
$ws.Range("R5").Value = 1.253975958396
$ws.Range("S5").Value = 0.07698529598963252
$ws.Range("T5").Value = 0.07698529598963252
$ws.Range("G6").Value = 0.015206
$ws.Range("H6").Value = 0.045618
$ws.Range("I6").Value = 0.197610548932631
$ws.Range("J6").Value = 0.197610548932631
$ws.Range("M6").Value = 8.050406666666667
$ws.Range("O6").Value = 0.3422817660053286
$ws.Range("P6").Value = 0.3422817660053286
$ws.Range("Q6").Value = 0.1224144837733333
$ws.Range("R6").Value = 1.10173035396
$ws.Range("S6").Value = 0.06763848766994332
$ws.Range("T6").Value = 0.06763848766994333
$ws.Range("G7").Value = 0.015206
$ws.Range("H7").Value = 0.045618
$ws.Range("I7").Value = 0.197610548932631
$ws.Range("J7").Value = 0.197610548932631
$ws.Range("M7").Value = 6.306542666666666
$ws.Range("N7").Value = 18.919628
$ws.Range("O7").Value = 0.2681373315304097
$ws.Range("P7").Value = 0.2681373315304097
$ws.Range("Q7").Value = 0.09589728778933332
$ws.Range("R7").Value = 0.8630755901039999
$ws.Range("S7").Value = 0.05298676527305512
$ws.Range("T7").Value = 0.05298676527305513
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.05034566666666667
$ws.Range("H8").Value = 0.151037
$ws.Range("I8").Value = 0.6542703423897976
$ws.Range("J8").Value = 0.6542703423897976
$ws.Range("M8").Value = 9.162874
$ws.Range("N8").Value = 27.488622
$ws.Range("O8").Value = 0.3895809024642617
$ws.Range("P8").Value = 0.3895809024642617
$ws.Range("Q8").Value = 0.4613110001126667
$ws.Range("R8").Value = 4.151799001014
$ws.Range("S8").Value = 0.2548912304438188
$ws.Range("T8").Value = 0.2548912304438188
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.05034566666666667
$ws.Range("H9").Value = 0.151037
$ws.Range("I9").Value = 0.6542703423897976
$ws.Range("J9").Value = 0.6542703423897976
$ws.Range("M9").Value = 8.050406666666667
$ws.Range("O9").Value = 0.3422817660053286
$ws.Range("P9").Value = 0.3422817660053286
$ws.Range("Q9").Value = 0.4053030905711112
$ws.Range("R9").Value = 3.647727815140001
$ws.Range("S9").Value = 0.2239448082380909
$ws.Range("T9").Value = 0.2239448082380909
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.05034566666666667
$ws.Range("H10").Value = 0.151037
$ws.Range("I10").Value = 0.6542703423897976
$ws.Range("J10").Value = 0.6542703423897976
$ws.Range("M10").Value = 6.306542666666666
$ws.Range("N10").Value = 18.919628
$ws.Range("O10").Value = 0.2681373315304097
$ws.Range("P10").Value = 0.2681373315304097
$ws.Range("Q10").Value = 0.3175070949151111
$ws.Range("R10").Value = 2.857563854236
$ws.Range("S10").Value = 0.1754343037078878
$ws.Range("T10").Value = 0.1754343037078878
